$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D header
$ws.Range("D1").Value = "Resultado"

# Mark existing rows' result column
$ws.Range("D2").Value = "x"
$ws.Range("D3").Value = "x"
$ws.Range("D4").Value = "x"

# New rows of data
$ws.Range("A5").Value = "TEST4"
$ws.Range("B5").Value = "hola"
$ws.Range("C5").Value = "mundo"
$ws.Range("D5").Value = "x"

$ws.Range("A6").Value = "TEST5"
$ws.Range("B6").Value = "tomas "
$ws.Range("C6").Value = "vino}"
$ws.Range("D6").Value = "x"

# B2 no longer carries its old explicit (non-visual) style override
$ws.Range("B2").Style = "Normal"

# Update active selection to D2
$ws.Range("D2").Select() | Out-Null
